$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The standalone row containing only "5840514 - Graziela Zamponi" (old row 13)
# was removed, shifting every subsequent row up by one.
$ws.Rows.Item(13).Delete()

# Content updates on the rows that remain (indices below are POST-delete,
# i.e. matching the final layout).
$ws.Range("B10").Value = "5840514 - Graziela Zamponi"
$ws.Range("C10").Value = "5840514 - Graziela Zamponi"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

$ws.Range("B18").Value = "5840514 - Graziela Zamponi"
$ws.Range("C18").Value = "5840514 - Graziela Zamponi"

$ws.Range("B21").Value = "Ao aluno que não alcançar a média 5,0 (cinco) no final do período letivo será dada uma nova oportunidade para a reelaboração dos trabalhos."
$ws.Range("C21").Value = "Ao aluno que não alcançar a média 5,0 (cinco) no final do período letivo será dada uma nova oportunidade para a reelaboração dos trabalhos."
